$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Sheet1"): row 1 changes from a 10-col numeric scratch row to
#     A1 = "Lion" (shared string), B1 = 3 (int-format style), rest cleared ---
$ws1.Range("A1").Value = "Lion"
$ws1.Range("B1").Value = 3
$ws1.Range("B1").NumberFormat = "0"
$ws1.Range("C1:J1").ClearContents()

# --- Sheet2 ("testCitizen"): rename the "ulkeleris.." / "umis.." strings ---
$ws2.Range("A1").Value = "ulkemis11"
$ws2.Range("A2").Value = "ulkemis22"
$ws2.Range("A3").Value = "ulkemis33"
$ws2.Range("A4").Value = "ulkemis44"
$ws2.Range("A5").Value = "ulkemis55"
$ws2.Range("A6").Value = "ulkemis66"
$ws2.Range("A7").Value = "ulkemis77"
$ws2.Range("A8").Value = "ulkemis88"

$ws2.Range("B1").Value = "uis11"
$ws2.Range("B2").Value = "uis21"
$ws2.Range("B3").Value = "uis31"
$ws2.Range("B4").Value = "uis41"
$ws2.Range("B5").Value = "uis51"
$ws2.Range("B6").Value = "uis61"
$ws2.Range("B7").Value = "uis71"
$ws2.Range("B8").Value = "uis81"

# --- Sheet2: widen column A slightly, and move the selection/active cell ---
$ws2.Columns.Item(1).ColumnWidth = 19.45
$ws2.Range("B9:B10").Select()
